# namc_customer_report.xlsx -- v.0.1.48 staging deploy
# Adds a "Director" contact block (Director / Telephone / Email) above the
# existing "Report Prepared For" section, fixes the "Telelphone:" typo
# further down, and shifts the dependent named ranges + hyperlinks
# accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Make room: insert 3 blank rows above the old spacer row 9. Every
#    row from the old row 9 onward (Report Prepared For, Customer info,
#    About this Report, etc.) shifts down by 3 rows, carrying its
#    existing values/styles with it.
# ---------------------------------------------------------------------
$ws.Rows("9:11").Insert()

# ---------------------------------------------------------------------
# 2. Populate the new Director contact block in rows 9-11.
#    The write order below matches the order new strings were appended
#    to the shared-string table by the original author.
# ---------------------------------------------------------------------
$ws.Range("D9").Value = "Director:"
$ws.Range("E9").Value = "Trip Armstrong"

$ws.Range("D10").Value = "Telephone:"

$ws.Range("D11").Value = "Email:"
$ws.Range("E11").Value = "trip.armstrong@usu.edu"
$ws.Range("E11").Style = "Hyperlink"

$ws.Range("E10").Value = "(760) 709-1210"

# Apply the existing right-aligned label style used throughout column D.
$ws.Range("D9,D10,D11").HorizontalAlignment = -4152

# ---------------------------------------------------------------------
# 3. Fix the long-standing "Telelphone:" typo (now living at D19 after
#    the row insertion above) so it reads "Telephone:".
# ---------------------------------------------------------------------
$ws.Range("D19").Value = "Telephone:"

# ---------------------------------------------------------------------
# 4. Hyperlinks did not auto-shift with the row insert, so rebuild the
#    hyperlink collection from scratch in the correct final positions
#    (preserving relationship order: report-link, namc-site, email).
# ---------------------------------------------------------------------
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A24"), "https://namc-usu.org/SampleProcessing/ResultsAndReports")
$ws.Hyperlinks.Add($ws.Range("D7"), "https://namc-usu.org/")
$ws.Hyperlinks.Add($ws.Range("E11"), "mailto:trip.armstrong@usu.edu")

# ---------------------------------------------------------------------
# 5. Update the defined names so they still point at the right cells
#    now that everything from row 9 down moved by +3 rows.
# ---------------------------------------------------------------------
$wb.Names.Item("SubmitterName").RefersTo = "=Summary!`$E`$14"
$wb.Names.Item("CustomerName").RefersTo  = "=Summary!`$E`$15"
$wb.Names.Item("Address1").RefersTo      = "=Summary!`$E`$16"
$wb.Names.Item("Address2").RefersTo      = "=Summary!`$E`$17"
$wb.Names.Item("CityStateZip").RefersTo  = "=Summary!`$E`$18"
$wb.Names.Item("Telephone").RefersTo     = "=Summary!`$E`$19"
$wb.Names.Item("Email").RefersTo         = "=Summary!`$E`$20"
$wb.Names.Item("ReportDate").RefersTo    = "=Summary!`$D`$12"

# ---------------------------------------------------------------------
# 6. Restore the active selection on the newly inserted Director cell.
# ---------------------------------------------------------------------
[void]$ws.Range("E13").Select()
